# Applies cryptos list update (prices/volumes refreshed, two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.561.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.304.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0911'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.967'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.33'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.652.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.309.97'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.493.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '281.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.79%  '
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +18.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("E28").Value = '  +3.85%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.137'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.38%  '
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0369'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.48%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("E39").Value = '  +2.20%  '
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.51%  '
$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -1.13%  '
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '112.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '78.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.610.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.55%  '
